$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 0.9600523903077931

# Add new row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 0.9590643274853801

# Add new row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 0.9393939393939394

# Add new row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 0.9691482226693494

# Copy the formatting from A2 (bold, bordered, centered style) down to A3:A5
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
